$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

# Columns A-D hold text values (date/time/weekday/week-number-as-text).
# Force text number format first so Excel does not auto-convert the
# date/time-looking strings (or the "06" code) into real dates/numbers.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-09"
$ws.Cells.Item($row, 2).Value = "11:42:24"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "06"

# Restore the default (unstyled) cell style so the new row matches the
# plain, style-less cells used by the rest of the data rows.
$textRange.Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126791
$ws.Cells.Item($row, 6).Value = 141916
$ws.Cells.Item($row, 7).Value = 168686
$ws.Cells.Item($row, 8).Value = 158392
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 143812
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191767
$ws.Cells.Item($row, 14).Value = 115100
$ws.Cells.Item($row, 15).Value = 44806
$ws.Cells.Item($row, 16).Value = 28411
$ws.Cells.Item($row, 17).Value = 64236
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42458
$ws.Cells.Item($row, 20).Value = -1
